# Fixed signin breaking, stylized home_frame, added cumulative earnings to statistics_frame
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert a new "Unnamed: 0.3" column before the existing "Unnamed: 0.2" column (old col B) ---
# This shifts old B..I to C..J, matching the dimension growing from I3 to J3.
$ws.Range("B1").EntireColumn.Insert()

# Give the new header cell (B1) the same style as the neighboring header cells before
# writing its text (copy format from C1, which carries the original header style s=1).
$ws.Range("C1").Copy()
$ws.Range("B1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("B1").Value2 = "Unnamed: 0.3"

# --- The two data rows (2 and 3) swap places: "Housing Purchase" moves to row 2 and
# "Direct Deposit" moves to row 3 (e.g. re-sorted by date). Swap via a scratch row so we
# don't lose the blank (empty-string) cells in columns D/E along the way. ---
$scratch = "A1000:J1000"
$ws.Range("A3:J3").Cut($ws.Range($scratch))
$ws.Range("A2:J2").Cut($ws.Range("A3:J3"))
$ws.Range($scratch).Cut($ws.Range("A2:J2"))
$ws.Range($scratch).EntireRow.Delete()

# --- New column B holds the cumulative-earnings statistic for each row. ---
# (The column insert copied column A's bold/bordered header style down into B2:B3;
# these are plain data cells like their neighbours, so drop back to the default style.)
$ws.Range("B2:B3").ClearFormats()
$ws.Range("B2").Value2 = 1
$ws.Range("B3").Value2 = 0

Write-Host "Edit applied."
